# Appends a new data row (row 32) to each of the four worksheets,
# mirroring the structure/style of the preceding row (row 31).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$newTime = 45818.46224537037

# Per-sheet row data: B, C, D, E (text), F, G (string form of the scientific
# literal so the exact double bit-pattern is reproduced), H, I (numbers)
$sheetsData = @(
    @{ Index = 1; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"; D = "0x01,0x84"; E = "0x07"; F = 400; G = "5.68631262647113e23"; H = 388; I = 7 },
    @{ Index = 2; B = "0x01,0x7c"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"; D = "0x01,0x74"; E = "0x19"; F = 380; G = "5.68432987514711e23"; H = 372; I = 25 },
    @{ Index = 3; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"; D = "0x00,0x6D"; E = "0x15"; F = 110; G = "5.68631262647113e23"; H = 109; I = 15 },
    @{ Index = 4; B = "0x00,0x82"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"; D = "0x00,0x81"; E = "0x9";  F = 130; G = "5.68631262647113e23"; H = 129; I = 9 }
)

foreach ($row in $sheetsData) {
    $ws = $wb.Worksheets.Item($row.Index)

    $gValue = $row.G -as [double]

    $ws.Range("A32").Value = $newTime
    $ws.Range("A32").NumberFormat = $dateFormat

    $ws.Range("B32").Value = $row.B
    $ws.Range("C32").Value = $row.C
    $ws.Range("D32").Value = $row.D
    $ws.Range("E32").Value = $row.E

    $ws.Range("F32").Value = $row.F
    $ws.Range("G32").Value = $gValue
    $ws.Range("H32").Value = $row.H
    $ws.Range("I32").Value = $row.I
}

Write-Host "Row 32 added to all sheets"
